# Apply updated NATMI Dll4-Notch4 values (per Dr Hou's advice)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns E..T (17 data columns) for rows 2..10.
# Row 1 is the header row and is left untouched.
$values = @(
    @(3, 1, 24.333402, 73.00020600000001, 0.9697175080062574, 0.9697175080062576, 3, 1, 30.345835, 91.03750500000001, 0.8527782452855476, 0.8527782452855475, 738.4174020806702, 6645.756618726031, 0.8269539949002501, 0.8269539949002501),
    @(3, 1, 24.333402, 73.00020600000001, 0.9697175080062574, 0.9697175080062576, 3, 1, 3.224107666666666, 9.672322999999999, 0.09060383010029813, 0.09060383010029811, 78.45350794428201, 706.081571498538, 0.08786012034068344, 0.08786012034068344),
    @(3, 1, 24.333402, 73.00020600000001, 0.9697175080062574, 0.9697175080062576, 3, 1, 2.014730333333333, 6.044191, 0.05661792461415433, 0.05661792461415433, 49.025243122594, 441.227188103346, 0.05490339276532388, 0.05490339276532389),
    @(3, 1, 0.288887, 0.866661, 0.01151252018667195, 0.01151252018667195, 3, 1, 30.345835, 91.03750500000001, 0.8527782452855476, 0.8527782452855475, 8.766517235645001, 78.898655120805, 0.009817626763604552, 0.009817626763604552),
    @(3, 1, 0.288887, 0.866661, 0.01151252018667195, 0.01151252018667195, 3, 1, 3.224107666666666, 9.672322999999999, 0.09060383010029813, 0.09060383010029811, 0.9314027915003332, 8.382625123502999, 0.001043078423019478, 0.001043078423019478),
    @(3, 1, 0.288887, 0.866661, 0.01151252018667195, 0.01151252018667195, 3, 1, 2.014730333333333, 6.044191, 0.05661792461415433, 0.05661792461415433, 0.5820294018056666, 5.238264616251, 0.0006518150000479225, 0.0006518150000479226),
    @(3, 1, 0.4710003333333333, 1.413001, 0.01876997180707065, 0.01876997180707065, 3, 1, 30.345835, 91.03750500000001, 0.8527782452855476, 0.8527782452855475, 14.29289840027833, 128.636085602505, 0.0160066236216929, 0.01600662362169291),
    @(3, 1, 0.4710003333333333, 1.413001, 0.01876997180707065, 0.01876997180707065, 3, 1, 3.224107666666666, 9.672322999999999, 0.09060383010029813, 0.09060383010029811, 1.518555785702555, 13.667002071323, 0.001700631336595215, 0.001700631336595215),
    @(3, 1, 0.4710003333333333, 1.413001, 0.01876997180707065, 0.01876997180707065, 3, 1, 2.014730333333333, 6.044191, 0.05661792461415433, 0.05661792461415433, 0.9489386585767776, 8.540447927191, 0.001062716848782528, 0.001062716848782528)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $rowValues = $values[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $j + 5  # column E is index 5
        $ws.Cells.Item($row, $col).Value2 = $rowValues[$j]
    }
}

$wb.Save()
